# Apply updated "F" column (ticket/interest count) values across the four
# worksheets, matching the regenerated site data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 33
$ws1.Range("F5").Value = 1722
$ws1.Range("F11").Value = 1172
$ws1.Range("F16").Value = 1311
$ws1.Range("F21").Value = 9
$ws1.Range("F24").Value = 13
$ws1.Range("F26").Value = 4301

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 6
$ws2.Range("F22").Value = 131
$ws2.Range("F46").Value = 29
$ws2.Range("F47").Value = 29

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F10").Value = 2974
$ws3.Range("F14").Value = 225

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 33
$ws4.Range("F9").Value = 2974
$ws4.Range("F15").Value = 1722
$ws4.Range("F27").Value = 6
$ws4.Range("F33").Value = 131
$ws4.Range("F34").Value = 131
$ws4.Range("F48").Value = 29
$ws4.Range("F49").Value = 29
